$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix up the header labels in row 1 (text tweaks) ---
$ws.Range("B1").Value = "X coords cm"
$ws.Range("C1").Value = "Y coords cm"
$ws.Range("F1").Value = "Weight gm"

# --- New "Jungle Bridge" vertex/weight data (rows 2-8) ---
$ws.Range("A2").Value = "Vertex 1"
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 0

$ws.Range("A3").Value = "Vertex 2"
$ws.Range("B3").Value = 7.2
$ws.Range("C3").Value = 5.5
$ws.Range("E3").Value = "Weight 1"
$ws.Range("F3").Value = 26

$ws.Range("A4").Value = "Vertex 3"
$ws.Range("B4").Value = 12
$ws.Range("C4").Value = 6.8
$ws.Range("E4").Value = "Weight 2"
$ws.Range("F4").Value = 31

$ws.Range("A5").Value = "Vertex 4"
$ws.Range("B5").Value = 20.9
$ws.Range("C5").Value = 7.5
$ws.Range("E5").Value = "Weight 3"
$ws.Range("F5").Value = 41

$ws.Range("A6").Value = "Vertex 5"
$ws.Range("B6").Value = 26.3
$ws.Range("C6").Value = 6.7
$ws.Range("E6").Value = "Weight 4"
$ws.Range("F6").Value = 46

$ws.Range("A7").Value = "Vertex 6"
$ws.Range("B7").Value = 29.2
$ws.Range("C7").Value = 3.3
$ws.Range("E7").Value = "Weight 5"
$ws.Range("F7").Value = 50

$ws.Range("A8").Value = "Vertex 7"
$ws.Range("B8").Value = 31.4
$ws.Range("C8").Value = 0

# B2 picks up the bold "header" cell style (s="1"), matching the original file.
$ws.Range("A1").Copy()
$ws.Range("B2").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Row heights: all data rows are 15.75pt (explicit/custom) ---
$ws.Rows.Item(2).RowHeight = 15.75
$ws.Rows.Item(3).RowHeight = 15.75
$ws.Rows.Item(4).RowHeight = 15.75
$ws.Rows.Item(5).RowHeight = 15.75
$ws.Rows.Item(6).RowHeight = 15.75
$ws.Rows.Item(7).RowHeight = 15.75
$ws.Rows.Item(8).RowHeight = 15.75

# --- Remove the (now unused) chart/drawing placeholder ---
for ($i = $ws.Shapes.Count; $i -ge 1; $i--) {
    $ws.Shapes.Item($i).Delete()
}

# --- Page margins all zeroed out, matching the committed layout ---
$ws.PageSetup.LeftMargin = 0
$ws.PageSetup.RightMargin = 0
$ws.PageSetup.TopMargin = 0
$ws.PageSetup.BottomMargin = 0
$ws.PageSetup.HeaderMargin = 0
$ws.PageSetup.FooterMargin = 0

# --- Selection left on F3, matching the saved cursor position ---
$ws.Range("F3").Select()
